# Diary update: add entries for 29 Jan 2020 (W) and 30 Jan 2020 (Th)
# (rows 23-27 of Sheet1, which were previously blank placeholder rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: 29 Jan 2020 (W), 1540-1710 -------------------------------
$ws.Range("A23").Value = "29 Jan 2020 (W)"
$ws.Range("B23").Value = "1540-1710"
$ws.Range("C23").Value = "Harry, Deon, Thuc"
$ws.Range("D23").Value = "Work on third lecture’s homework and write up the component findings"
$ws.Range("E23").Value = "Deon printed the UML diagram and taped everything together. It’s as messy as we expected. We chose which components to highlight. We wrote up on where the components fit in the UML diagram."
$ws.Range("G23").Value = "I came in 30 minutes early. Teammates came 40 minutes late. People are busy and traffic is not nice. The difference in time gives a few moments to breathe and observe the lab."
$ws.Rows.Item(23).RowHeight = 73.1

# --- Row 24: 29 Jan 2020 (W), 1910-1945 -------------------------------
$ws.Range("A24").Value = "29 Jan 2020 (W)"
$ws.Range("B24").Value = "1910-1945"
$ws.Range("C24").Value = "Harry, Deon, Thuc"
$ws.Range("D24").Value = "We continued discussing the diagram and what components to select"
$ws.Range("E24").Value = "We chose Metronome* and ScreenMarkers and found Metronome* on the diagram"
$ws.Rows.Item(24).RowHeight = 37.3

# --- Row 25: 29 Jan 2020 (W), 2300-0000 -------------------------------
$ws.Range("A25").Value = "29 Jan 2020 (W)"
$ws.Range("B25").Value = "2300-0000"
$ws.Range("C25").Value = "Harry"
$ws.Range("D25").Value = "Finalized the writeup"
$ws.Range("E25").Value = "Finalized the writeup"

# --- Row 26: 30 Jan 2020 (Th), 0000-0200 ------------------------------
$ws.Range("A26").Value = "30 Jan 2020 (Th)"
$ws.Range("B26").Value = "0000-0200"
$ws.Range("C26").Value = "Harry"
$ws.Range("D26").Formula = "=D25"
$ws.Range("E26").Formula = "=E25"
$ws.Range("G26").Value = "More tired than usual. Nutrition and exercise aren’t the issue; something else is wrong."
$ws.Rows.Item(26).RowHeight = 37.3

# --- Row 27: 30 Jan 2020 (Th), 1135-1220 ------------------------------
$ws.Range("A27").Value = "30 Jan 2020 (Th)"
$ws.Range("B27").Value = "1135-1220"
$ws.Range("C27").Value = "Harry"
$ws.Range("D27").Value = "Look at the homework once more. Will look at UML diagram later today to make sure everything is correct."
$ws.Range("E27").Value = "Things look good to go."
$ws.Rows.Item(27).RowHeight = 37.3

# --- Update view: scroll down a bit, select G27 -----------------------
$ws.Range("G27").Select() | Out-Null
